# Update sheet name and data to reflect new "through" date of 2022-08-05
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name), reflecting the new "through" date
$ws.Name = "Through 2022-07-28"

# Update the header label in I1 ("2022 (through 07-27)" -> "2022 (through 07-28)")
$ws.Range("I1").Value = "2022 (through 07-28)"

# Update the data values for the new day's data
$ws.Range("I8").Value = 158
$ws.Range("I14").Value = 964
